$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (meanEMG) values for columns B:E
$ws.Range("B2").Value = 12.600043675037883
$ws.Range("C2").Value = 12.02188575120635
$ws.Range("D2").Value = 13.362599958526902
$ws.Range("E2").Value = 12.944287081014441

# Update row 3 (legmaxROM) values for columns B:E
$ws.Range("B3").Value = 12.855681704249111
$ws.Range("C3").Value = 10.818102188479116
$ws.Range("D3").Value = 15.074578803049086
$ws.Range("E3").Value = 12.20806682471637

# Update the selection to reflect the narrower highlighted range
$ws.Range("B1:E3").Select()
